# 3rd commit by Pradnya
# Updates test AWB serial numbers / status codes on the "Add_AWB" and
# "View_AWB" sheets, removes an obsolete test row from "Add_AWB", and
# tweaks column widths / selections left over from manual editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Add_AWB": refresh serial numbers, a couple of field values, and
# drop the now-unused 6th test row.
# ---------------------------------------------------------------------
$wsAdd = $wb.Worksheets.Item("Add_AWB")

$wsAdd.Range("B2").Value = "'9702484"

$wsAdd.Range("B3").Value = "'9702485"
$wsAdd.Range("E3").Value = "CC"

$wsAdd.Range("B4").Value = "'9702486"
$wsAdd.Range("J4").Value = "'2121212121"

$wsAdd.Range("B5").Value = "'9702487"
$wsAdd.Range("J5").Value = ""

$wsAdd.Rows("6:6").Delete()

$wsAdd.Columns("C").ColumnWidth = 13

# ---------------------------------------------------------------------
# Sheet "View_AWB": refresh serial numbers and flown/export billing
# status codes for rows 3 and 4.
# ---------------------------------------------------------------------
$wsView = $wb.Worksheets.Item("View_AWB")

$wsView.Range("B3").Value = "'9702479"

$wsView.Range("B4").Value = "'9702479"
$wsView.Range("C4").Value = "CC"
$wsView.Range("D4").Value = "N"
$wsView.Range("E4").Value = "NC"

# ---------------------------------------------------------------------
# Restore the on-disk selections (View_AWB selected first so Add_AWB
# ends up as the active/tab-selected sheet, matching the saved file).
# ---------------------------------------------------------------------
$wsView.Range("C9").Select()
$wsAdd.Range("C10").Select()
